$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg).
# These values correspond to a weekly re-shuffle of the existing rows'
# date-linked data (Fruta / hortaliza, semanal).
$rowData = @{
    2  = @{ D = 45096; M = 30; N = 20000; O = 20000; P = 20000; S = 1111 }
    3  = @{ D = 45083; M = 50; N = 15000; O = 15000; P = 15000; S = 833 }
    4  = @{ D = 45069; M = 60; N = 15000; O = 15000; P = 15000; S = 833 }
    5  = @{ D = 45084; M = 50; N = 18000; O = 19000; P = 18500; S = 1028 }
    6  = @{ D = 45076; M = 20; N = 15000; O = 15000; P = 15000; S = 833 }
    7  = @{ D = 45061; M = 40; N = 15000; O = 15000; P = 15000; S = 833 }
    8  = @{ D = 45092; M = 60; N = 18000; O = 19000; P = 18667; S = 1037 }
    9  = @{ D = 45085; M = 30; N = 19000; O = 19000; P = 19000; S = 1056 }
    10 = @{ D = 45055; M = 50; N = 15000; O = 15000; P = 15000; S = 833 }
    11 = @{ D = 45111; M = 20; N = 20000; O = 20000; P = 20000; S = 1111 }
    12 = @{ D = 45112; M = 20; N = 20000; O = 20000; P = 20000; S = 1111 }
    13 = @{ D = 45072; M = 30; N = 15000; O = 15000; P = 15000; S = 833 }
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
